$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") for rows 2-5 moves from serial 45207 (2023-10-08)
# to serial 45208 (2023-10-09). The cells already carry a date number
# format, so assigning a Date value keeps that formatting intact.
foreach ($row in 2..5) {
    $ws.Cells.Item($row, 3).Value = 45208
}
